# Updated cryptos list with latest coinranking.com snapshot (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.840.01"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.637.55"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'215.49"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").Value = "'0.5056"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.06422"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").Value = "'19.69"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").Value = "'0.07774"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").Value = "'4.282"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "1.639.08"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value = "1.862.39"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "'0.5601"
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("D16").Value = "0.0₅7597"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").Value = "'63.01"
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("D18").Value = "25.837.62"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "'195.12"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'4.316"
$ws.Range("E21").Value = "  -2.33%  "
$ws.Range("D22").Value = "'9.878"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "'6.106"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'1.780"
$ws.Range("E25").Value = "  -6.35%  "
$ws.Range("D26").Value = "'140.16"
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("D27").Value = "'0.1265"
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("D28").Value = "'6.808"
$ws.Range("D29").Value = "'15.43"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D31").Value = "'0.04871"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("D32").Value = "'3.295"
$ws.Range("E32").Value = "  +2.00%  "
$ws.Range("D33").Value = "'3.223"
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("D35").Value = "'2.375"
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").Value = "'0.9020"
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("D37").Value = "'2.576"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.127.48"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'0.5509"
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").Value = "'0.9953"
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("D42").Value = "'5.535"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "'0.8021"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "'97.89"
$ws.Range("D45").Value = "1.773.37"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").Value = "0.0₈112"
$ws.Range("E46").Value = "  -9.06%  "
$ws.Range("D47").Value = "'55.30"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("D48").Value = "'0.4279"
$ws.Range("E48").Value = "  -4.35%  "
$ws.Range("D49").Value = "'7.684"
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").Value = "'0.05045"
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "  +0.05%  "
